$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 248, pushing existing rows 248..303 down to 249..304
$ws.Rows.Item(248).Insert()

# Populate the newly inserted row 248 with the new record
$ws.Cells.Item(248, 1).Value = 9
$ws.Cells.Item(248, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(248, 3).Value = "Metropolitana"
$ws.Cells.Item(248, 4).Value = 44736
$ws.Cells.Item(248, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(248, 5).Value = 13
$ws.Cells.Item(248, 6).Value = 100112021
$ws.Cells.Item(248, 7).Value = "Ají"
$ws.Cells.Item(248, 8).Value = "Americana (o)"
$ws.Cells.Item(248, 9).Value = "Primera"
$ws.Cells.Item(248, 10).Value = 25
$ws.Cells.Item(248, 11).Value = 39000
$ws.Cells.Item(248, 12).Value = 40000
$ws.Cells.Item(248, 13).Value = 39480
$ws.Cells.Item(248, 14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(248, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(248, 16).Value = 1579
$ws.Cells.Item(248, 17).Value = 25
$ws.Cells.Item(248, 18).Value = "Hortaliza"
